$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11; this pushes the existing rows 11-63 down to 12-64,
# and copies formatting from the row above (row 10) by default.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with a fresh data record.
# Columns that are unchanged vs. the row that used to occupy row 11
# (now shifted to row 12) keep the same values; D, J, K, L, M, P get new values.
$ws.Cells.Item(11, 1).Value = 4
$ws.Cells.Item(11, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(11, 3).Value = "Los Lagos"
$ws.Cells.Item(11, 4).Value = (Get-Date -Year 2023 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 100112030
$ws.Cells.Item(11, 7).Value = "Poroto granado"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 35000
$ws.Cells.Item(11, 12).Value = 35000
$ws.Cells.Item(11, 13).Value = 35000
$ws.Cells.Item(11, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(11, 15).Value = "Región Metropolitana"
$ws.Cells.Item(11, 16).Value = 1400
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = "Hortaliza"
